# Data for Development.xlsx - apply the commit's edits via Excel COM-interop.
#
# Summary of the change (per the commit message / xml diff):
#  - The "Bullet Message" sheet is removed entirely (bullet message protocol
#    dropped in favor of the new game-status / score tracking work).
#  - On "Object Code": the "bullet" object-code row now represents the new
#    "game status" message, and the now-unused "AI" / "Item draw" rows are
#    cleared out.
#  - Selection/active-tab bookkeeping moves from "Weapons" (A10 selected) to
#    "Object Code" (now the active tab, C3 selected); "Weapons" keeps A2
#    selected instead.

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove the obsolete "Bullet Message" sheet ---------------------------
$bulletSheet = $wb.Worksheets.Item("Bullet Message")
$bulletSheet.Delete()

# --- Rework the "Object Code" table ----------------------------------------
$objectCode = $wb.Worksheets.Item("Object Code")
$objectCode.Range("B3").Value = "game status"   # was "bullet"
$objectCode.Range("B5").ClearContents()          # was "AI"
$objectCode.Range("B6").ClearContents()          # was "Item draw"

# --- Update selections / active sheet --------------------------------------
$weapons = $wb.Worksheets.Item("Weapons")
$weapons.Activate()
$weapons.Range("A2").Select()

$objectCode.Activate()
$objectCode.Range("C3").Select()
